$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row for "2035_TM152_FBP_Plus_14" right after the existing
#    2035 FBP Plus block (old row 61), re-using row 60's formatting so the
#    "last row of block" border style (previously on row 61) slides down
#    onto the new row 62, exactly like Excel's own Insert-Copied-Cells flow.
# ---------------------------------------------------------------------------
$ws.Rows.Item(60).Copy()
$ws.Rows.Item(61).Insert()

# Row 61 now holds a duplicate of old row 60; restore it to the original
# "2035_TM152_FBP_Plus_13" content that used to live there.
$ws.Range("A61").Value = "RTP2021"
$ws.Range("B61").Value = 2035
$ws.Range("C61").Value = "2035_TM152_FBP_Plus_13"
$ws.Range("D61").Value = "FinalBlueprint"
$ws.Range("E61").Value = "Plus"
$ws.Range("F61").Value = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.22"'
$ws.Range("G61").Value = "run352"
$ws.Range("H61").Value = "current"

# Row 62 (pushed down from the old row 61) kept its original border/style but
# still has the old "_13" values -- overwrite with the new "_14" series.
$ws.Range("A62").Value = "RTP2021"
$ws.Range("B62").Value = 2035
$ws.Range("C62").Value = "2035_TM152_FBP_Plus_14"
$ws.Range("D62").Value = "FinalBlueprint"
$ws.Range("E62").Value = "Plus"
$ws.Range("F62").Value = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25"'
$ws.Range("G62").Value = "run182"
$ws.Range("H62").Value = "current"

# ---------------------------------------------------------------------------
# 2) Append a new row for "2050_TM152_PlusCrossing_14" at the very bottom of
#    the table (new row 98). Use copy+Insert (not a plain Paste) so the
#    formatting of the last existing row carries over correctly even though
#    there is nothing below it to shift.
# ---------------------------------------------------------------------------
$ws.Rows.Item(96).Copy()
$ws.Rows.Item(97).Insert()

$ws.Rows.Item(96).Copy()
$ws.Rows.Item(98).Insert()

# Row 97 is a duplicate of old row 96 ("2050_TM152_FBP_PlusCrossing_13") --
# that is exactly the content that belongs there, so no edits are needed.

# Row 98 is also a duplicate of old row 96 -- overwrite with the new "_14"
# series values.
$ws.Range("A98").Value = "RTP2021"
$ws.Range("B98").Value = 2050
$ws.Range("C98").Value = "2050_TM152_FBP_PlusCrossing_14"
$ws.Range("D98").Value = "FinalBlueprint"
$ws.Range("E98").Value = "Plus"
$ws.Range("F98").Value = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25"'
$ws.Range("G98").Value = "run182"
$ws.Range("H98").Value = "current"

# Match the source workbook's small end-of-table bottom border nuance on the
# newly appended last row's G cell (style carried a thin bottom border there).
$ws.Range("G98").Borders.Item(9).LineStyle = 1
